# Pivots workbook: add the new trading day's data into column J
# (the workbook already has columns E..I populated for prior days; this
# mirrors the previous "today" column, I, into the new "today" column, J
# -- same as a user copying I2:I46 -> J2:J46 and then typing in the new
# day's High/Low/Close plus the fixed EW labels).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Povit")

# 1) Copy formatting (fills/borders/number formats) from column I to column J
#    for the two blocks of rows that are used on this sheet (row 5 has no
#    G:J cells at all, so it is skipped).
$ws.Range("I2:I4").Copy()
$ws.Range("J2:J4").PasteSpecial(-4122)
$ws.Range("I6:I46").Copy()
$ws.Range("J6:J46").PasteSpecial(-4122)

# 2) New day's raw inputs (High / Low / Close)
$ws.Range("J2").Value = 10540.65
$ws.Range("J3").Value = 10261.9
$ws.Range("J4").Value = 10316.450000000001

# 3) Pivot-point formulas, mirrored from column I into column J
$ws.Range("J6").Formula = "=J8+J43"
$ws.Range("J7").Formula = "=J11+J43"
$ws.Range("J8").Formula = "=(2*J11)-J3"

$ws.Range("J10").Formula = "=J11+J13/2"
$ws.Range("J11").Formula = "=(J2+J3+J4)/3"
$ws.Range("J12").Formula = "=J11-J13/2"
$ws.Range("J13").Formula = "=ABS((J11-J46)*2)"

$ws.Range("J15").Formula = "=2*J11-J2"
$ws.Range("J16").Formula = "=J11-J43"
$ws.Range("J17").Formula = "=J15-J43"

# Camarilla pivots
$ws.Range("J19").Formula = "=(J2/J3)*J4"
$ws.Range("J20").Formula = "=J21+1.168*(J21-J22)"
$ws.Range("J21").Formula = "=J4+J44/2"
$ws.Range("J22").Formula = "=J4+J44/4"
$ws.Range("J23").Formula = "=J4+J44/6"
$ws.Range("J24").Formula = "=J4+J44/12"
$ws.Range("J25").Formula = "=J4"
$ws.Range("J26").Formula = "=J4-J44/12"
$ws.Range("J27").Formula = "=J4-J44/6"
$ws.Range("J28").Formula = "=J4-J44/4"
$ws.Range("J29").Formula = "=J4-J44/2"
$ws.Range("J30").Formula = "=J29-1.168*(J28-J29)"
$ws.Range("J31").Formula = "=J4-(J19-J4)"

# Elliott - Fibonacci "Close" row (this one already existed as a shared
# formula covering H37:N37; give J37 its own explicit formula)
$ws.Range("J37").Formula = "=J4"

# High-Low helper block
$ws.Range("J43").Formula = "=ABS(J2-J3)"
$ws.Range("J44").Formula = "=J43*1.1"
$ws.Range("J45").Formula = "=(J2+J3)"
$ws.Range("J46").Formula = "=(J2+J3)/2"

# 4) Selection / view bookkeeping: the sheet had scrolled to show row 23 at
#    the top with I41 selected; it is now back at the default scroll
#    position with M31 selected.
$ws.Range("M31").Select()
